$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the raw-data year rows (order matters: "2019(s)" must be interned
# into the shared-strings table before "2018(t)" to match the target indices)
$ws.Range("A5").Value = "2019(s)"
$ws.Range("A4").Value = "2018(t)"

# Fix the linked raw data for 2019 (row 5): c(mango) and c(beer)
$ws.Range("C5").Value = 20
$ws.Range("E5").Value = 50

# Reflect the author's view state at save time (selection + zoom)
$ws.Range("C6").Select() | Out-Null
$excel.ActiveWindow.Zoom = 115
